$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "Objetivos:" (row 10) value - it had the wrong text (a leftover
#    professor name) and should hold the actual Portuguese objectives text.
$ws.Range("B10").Value = "Propiciar conhecimentos básicos sobre os materiais terrestres e os principais processos geológicos."
$ws.Range("C10").Value = "Propiciar conhecimentos básicos sobre os materiais terrestres e os principais processos geológicos."

# 2) Insert a new row at 13 (pushes old rows 13-21 down to 14-22) to hold the
#    "Docentes responsáveis:" value (professor name), which previously lived
#    (mis-placed) in row 10/18.
$ws.Rows("13").Insert()

$ws.Range("B13").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C13").Value = "5464150 - Mariana Consiglio Kasemodel"
# match formatting of the other "value" columns (style 2 / style 3)
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) "Programa resumido:" (now row 14) gets its real short-syllabus text
#    instead of the placeholder "Semestral".
$ws.Range("B14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."
$ws.Range("C14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."

# 4) "Programa:" (now row 16) gets its real full-syllabus text instead of the
#    placeholder date.
$ws.Range("B16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."
$ws.Range("C16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."

# 5) "Método:" (now row 19) gets its real method text instead of the
#    misplaced professor name.
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# 6) "Bibliografia:" (now row 22) gets the real bibliography text instead of
#    the placeholder "1 (uma) prova escrita".
$ws.Range("B22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."
$ws.Range("C22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."

# Row-height touch-ups to mirror the published sheet layout.
$ws.Rows("13").RowHeight = 15
$ws.Rows("14").RowHeight = 60
$ws.Range("A14").RowHeight = 60
$ws.Rows("15").RowHeight = 60
$ws.Rows("16").RowHeight = 120
$ws.Rows("17").RowHeight = 120
$ws.Rows("21").RowHeight = 60
